$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 345, pushing the existing rows 345-360 down to
# 346-361 (this also mirrors the "Precio Fruta Hortalizas" weekly refresh
# pattern where a new week's record is prepended to this variety's block).
$ws.Rows.Item(345).EntireRow.Insert()

# Populate the newly inserted row 345 with this week's record.
$ws.Cells.Item(345, 1).Value = 4
$ws.Cells.Item(345, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(345, 3).Value = "Los Lagos"
$ws.Cells.Item(345, 4).Value = 44939
$ws.Cells.Item(345, 5).Value = 10
$ws.Cells.Item(345, 6).Value = 100112021
$ws.Cells.Item(345, 7).Value = "Ají"
$ws.Cells.Item(345, 8).Value = "Inferno"
$ws.Cells.Item(345, 9).Value = "Primera"
$ws.Cells.Item(345, 10).Value = 180
$ws.Cells.Item(345, 11).Value = 20000
$ws.Cells.Item(345, 12).Value = 20000
$ws.Cells.Item(345, 13).Value = 20000
$ws.Cells.Item(345, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(345, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(345, 16).Value = 2000
$ws.Cells.Item(345, 17).Value = 10
$ws.Cells.Item(345, 18).Value = "Hortaliza"
